# Apply the commit: add new "sigma_050" results sheet (the NLM-LBP run for
# sigma=0.50 is now complete) and refresh the higher-precision PSNR values
# in the existing "sigma_010" / "sigma_025" sheets.

$wb = $excel.ActiveWorkbook

# --- Refresh "sigma_010" (Sheet2) Noisy / NLM-LBP columns ---
$ws2 = $wb.Worksheets.Item("sigma_010")
$ws2.Cells.Item(2, 2).Value = 27.86556592171642
$ws2.Cells.Item(2, 3).Value = 28.88167433807014
$ws2.Cells.Item(3, 2).Value = 27.89314816732099
$ws2.Cells.Item(3, 3).Value = 28.88544476473718
$ws2.Cells.Item(4, 2).Value = 27.88777339169984
$ws2.Cells.Item(4, 3).Value = 28.90342015758691
$ws2.Cells.Item(5, 2).Value = 27.87356436109086
$ws2.Cells.Item(5, 3).Value = 28.85832136261758
$ws2.Cells.Item(6, 2).Value = 27.87101380393653
$ws2.Cells.Item(6, 3).Value = 28.88077634761183
$ws2.Cells.Item(7, 2).Value = 27.84684309650138
$ws2.Cells.Item(7, 3).Value = 28.87286763784935
$ws2.Cells.Item(8, 2).Value = 27.86665875623024
$ws2.Cells.Item(8, 3).Value = 28.89158522516483
$ws2.Cells.Item(9, 2).Value = 27.85289657898602
$ws2.Cells.Item(9, 3).Value = 28.91444071076244
$ws2.Cells.Item(10, 2).Value = 27.86844424876679
$ws2.Cells.Item(10, 3).Value = 28.88450866179215
$ws2.Cells.Item(11, 2).Value = 27.86511721561606
$ws2.Cells.Item(11, 3).Value = 28.87648556293106
$ws2.Cells.Item(12, 2).Value = 27.86910255418651
$ws2.Cells.Item(12, 3).Value = 28.88495247691235

# --- Refresh "sigma_025" (Sheet3) Noisy / NLM-LBP columns ---
$ws3 = $wb.Worksheets.Item("sigma_025")
$ws3.Cells.Item(2, 2).Value = 19.76110814097889
$ws3.Cells.Item(2, 3).Value = 25.15513782666984
$ws3.Cells.Item(3, 2).Value = 19.76119084291836
$ws3.Cells.Item(3, 3).Value = 25.1880382976403
$ws3.Cells.Item(4, 2).Value = 19.76618412944767
$ws3.Cells.Item(4, 3).Value = 25.07569095645501
$ws3.Cells.Item(5, 2).Value = 19.73913588460996
$ws3.Cells.Item(5, 3).Value = 25.07451578196252
$ws3.Cells.Item(6, 2).Value = 19.75397571247449
$ws3.Cells.Item(6, 3).Value = 25.06910352329858
$ws3.Cells.Item(7, 2).Value = 19.75419999468854
$ws3.Cells.Item(7, 3).Value = 25.08150014794403
$ws3.Cells.Item(8, 2).Value = 19.74023001408807
$ws3.Cells.Item(8, 3).Value = 25.06840370953961
$ws3.Cells.Item(9, 2).Value = 19.75200249952362
$ws3.Cells.Item(9, 3).Value = 25.09617719895786
$ws3.Cells.Item(10, 2).Value = 19.75256481012362
$ws3.Cells.Item(10, 3).Value = 25.05323524850856
$ws3.Cells.Item(11, 2).Value = 19.75863365551985
$ws3.Cells.Item(11, 3).Value = 25.07353832102982
$ws3.Cells.Item(12, 2).Value = 19.75392256843731
$ws3.Cells.Item(12, 3).Value = 25.09353410120061

# --- Add the new "sigma_050" sheet as the last tab, mirroring the layout
#     of the other sigma_* sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "sigma_050"

$ws4.Cells.Item(1, 1).Value = "Rows"
$ws4.Cells.Item(1, 2).Value = "Noisy"
$ws4.Cells.Item(1, 3).Value = "NLM-LBP"
$ws4.Cells.Item(2, 1).Value = 0
$ws4.Cells.Item(2, 2).Value = 14.78599328414847
$ws4.Cells.Item(2, 3).Value = 20.31091138274269
$ws4.Cells.Item(3, 1).Value = 1
$ws4.Cells.Item(3, 2).Value = 14.7614917969002
$ws4.Cells.Item(3, 3).Value = 20.26855433782312
$ws4.Cells.Item(4, 1).Value = 2
$ws4.Cells.Item(4, 2).Value = 14.75489515597404
$ws4.Cells.Item(4, 3).Value = 20.30142652972216
$ws4.Cells.Item(5, 1).Value = 3
$ws4.Cells.Item(5, 2).Value = 14.78516269221445
$ws4.Cells.Item(5, 3).Value = 20.28609537042992
$ws4.Cells.Item(6, 1).Value = 4
$ws4.Cells.Item(6, 2).Value = 14.77403269556762
$ws4.Cells.Item(6, 3).Value = 20.32449719411572
$ws4.Cells.Item(7, 1).Value = 5
$ws4.Cells.Item(7, 2).Value = 14.77679663735373
$ws4.Cells.Item(7, 3).Value = 20.2978428933309
$ws4.Cells.Item(8, 1).Value = 6
$ws4.Cells.Item(8, 2).Value = 14.76550662684715
$ws4.Cells.Item(8, 3).Value = 20.2867802760757
$ws4.Cells.Item(9, 1).Value = 7
$ws4.Cells.Item(9, 2).Value = 14.76390061800428
$ws4.Cells.Item(9, 3).Value = 20.30813913331753
$ws4.Cells.Item(10, 1).Value = 8
$ws4.Cells.Item(10, 2).Value = 14.76769493689635
$ws4.Cells.Item(10, 3).Value = 20.29789265260527
$ws4.Cells.Item(11, 1).Value = 9
$ws4.Cells.Item(11, 2).Value = 14.78641939021353
$ws4.Cells.Item(11, 3).Value = 20.28383041469881
$ws4.Cells.Item(12, 1).Value = "Média"
$ws4.Cells.Item(12, 2).Value = 14.77218938341198
$ws4.Cells.Item(12, 3).Value = 20.29659701848618
